$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace "... was rather low with most responding neutrally. "
#    with   "... was mixed between neutral and agreeing. "
# ------------------------------------------------------------------
$ok = $d.Content.Find.Execute(" was rather low with most responding neutrally. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " was mixed between neutral and agreeing. ", 2)

# Locate the freshly-written sentence so later searches stay inside it.
$outer = $d.Content
$ok = $outer.Find.Execute(" was mixed between neutral and agreeing. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$outerStart = $outer.Start
$outerEnd = $outer.End

# ------------------------------------------------------------------
# Split the single run produced by the replace above into the four
# separate runs Word creates when text is typed/edited incrementally:
#   " was " | "mixed between" | " " | "neutral and agreeing" | ". "
# Nudging a formatting property off/on forces a run boundary at the
# edges of the found text without altering the visible formatting.
# ------------------------------------------------------------------
$cursor = $outerStart

$r1 = $d.Range($cursor, $outerEnd)
$ok = $r1.Find.Execute("mixed between", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Bold = 1
$r1.Bold = 0
$cursor = $r1.End

$r2 = $d.Range($cursor, $outerEnd)
$ok = $r2.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Bold = 1
$r2.Bold = 0
$cursor = $r2.End

$r3 = $d.Range($cursor, $outerEnd)
$ok = $r3.Find.Execute("neutral and agreeing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Bold = 1
$r3.Bold = 0
$cursor = $r3.End

$r4 = $d.Range($cursor, $outerEnd)
$ok = $r4.Find.Execute(". ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4.Bold = 1
$r4.Bold = 0

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark (Word's last-edit marker) so it now
#    sits right after "neutral and agreeing" -- where the new text
#    was typed -- instead of its old spot after "... stall quite
#    frequently." Remove the stale bookmark, then recreate it at the
#    new location (re-adding a bookmark with the same name also
#    relocates it, but we delete explicitly first to be safe).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldGoBack = $d.Bookmarks("_GoBack")
    $oldGoBack.Delete()
}
$bmRange = $d.Range($r3.End, $r3.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
